$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.641.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.450.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.15%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.31%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.452.21"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.26%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.80"
$ws.Range("D10").Style = "Normal"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.128"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.392"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.040.48"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.10%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.51%  "

$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000176"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.41%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.451.02"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.05%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.770.51"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.16%  "

$ws.Range("E19").Value = "  +8.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.42"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.01%  "

$ws.Range("E21").Value = "  +0.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "390.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.13%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.567"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.587.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.81%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "73.19"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.86%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  +0.81%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000126"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.181"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.46%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.54"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -12.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("E33").Value = "  +1.44%  "

$ws.Range("E34").Value = "  +1.29%  "

$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "24.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.05"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.42%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.52%  "

$ws.Range("E39").Value = "  +1.54%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "166.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0789"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.32%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +11.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.794"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.09%  "

$ws.Range("E44").Value = "  +2.34%  "

$ws.Range("E45").Value = "  -0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "42.38"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "

$ws.Range("E47").Value = "  +0.61%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.606.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.49%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.26%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.14%  "
